$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the date values in column A (rows 2 and 3) with the plain text "2018".
# Forcing the number format to Text ("@") first stops Excel from re-interpreting
# the digit string "2018" back into a number; resetting the style afterwards
# clears the now-unneeded formatting so the cell goes back to the default look.
$rng = $ws.Range("A2:A3")
$rng.NumberFormat = "@"
$ws.Range("A2").Value = "2018"
$ws.Range("A3").Value = "2018"
$rng.Style = "표준"
